# Applies the edits described by the commit diff to the active document.
#
# Summary of changes:
#  1. Fix a couple of small wording mistakes.
#  2. Merge several runs that only differed by run-splitting (no textual
#     change) back into single runs.
#  3. Insert a new "Lame basale" paragraph (with a partially character
#     styled run) right after the Apicale/Basale table.
#  4. Move the <w:lastRenderedPageBreak/> marker from the "Sanguin" run to
#     the "Cartilageux" run.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $r = $d.Content
    $ok = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $oldText"
    }
    return $ok
}

# --- 1. Wording fixes -------------------------------------------------
Replace-Text " étude des structure des biologique." " étude des structures biologiques."
Replace-Text "Les formes des organes sont participe directement à leur fonction." "Les formes des organes sont adaptées à leur fonction."

# --- 2. Pure run-merges (old text == new text, just defragment runs) --
$mergeTexts = @(
    "La diversité des plans d’organisation corporelle est limitée par les lois de la physique. Dans certains cas, les contraintes du milieu ont imposé à la sélection naturelle a adopté la même la forme comme c’est le cas pour l’apparence fusiforme des animaux aquatiques.",
    "Les tissus épithéliaux sont constitués de cellules jointives serré. Elles forment des barrière. L’épithélium est polarisé, une face :",
    "couche unique transport de substances par diffusion appelé communément",
    "Des muqueuses, les voies respiratoire",
    "stratifié squameux",
    "protéger et maintenir les organes.",
    "Dont fait partie le tissu aréolaire épithélium aux tissus sous",
    "Contient un mélange de fibres de collagènes et chondroïtine sulfate qui confère à la fois résistance et souplesse",
    "La régulation se fait autour d’une valeur de référence qui peut varier au cours du temps notamment durant le rythme circadien."
)
foreach ($t in $mergeTexts) {
    Replace-Text $t $t
}

# --- 3. Insert the new "Lame basale" paragraph -------------------------
$nbsp = [char]0x00A0
$anchor = "Il existe plusieurs types de tissus épithéliaux" + $nbsp + ":"
$r = $d.Content
$found = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if (-not $found) {
    Write-Output "NOT FOUND: $anchor"
}
$r.InsertParagraphBefore()

$r2 = $d.Content
$r2.Find.Execute($anchor, $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$anchorPara = $r2.Paragraphs(1)
$newPara = $anchorPara.Previous()
$newRange = $newPara.Range
$newRange.Collapse(1)
$newText = "Lame basale assemblage de protéines et glycoprotéines qui permet l’adhérence au tissu conjonctif."
$newRange.InsertAfter($newText)

# Apply the "Accentuation" character style to the "Lame basale" label only,
# via Find/Replacement so no stray paragraph rsid gets stamped.
$label = "Lame basale"
$r3 = $d.Content
$r3.Find.ClearFormatting()
$r3.Find.Replacement.ClearFormatting()
$r3.Find.Replacement.Style = "Accentuation"
$r3.Find.Execute($label, $true, $false, $false, $false, $false, `
                  $true, 1, $true, $label, 2) | Out-Null

# --- 4. Move <w:lastRenderedPageBreak/> from "Sanguin" to "Cartilageux" -
$r4 = $d.Content
$r4.Find.Execute("Sanguin", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$sanguinPara = $r4.Paragraphs(1)
$sanguinText = $sanguinPara.Range.Text
$sanguinPara.Range.Text = $sanguinText

$r5 = $d.Content
$r5.Find.Execute("Cartilageux", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$cartPara = $r5.Paragraphs(1)
$cartRange = $cartPara.Range
$cartRange.Collapse(1)
$cartRange.Fields.Add($cartRange, 1) | Out-Null

Write-Output "done"
